$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.731.76'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.120.41'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''243.56'
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').Value = '''617.79'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').Value = '  -5.70%  '
$ws.Range('D8').Value = '''0.388'
$ws.Range('E8').Value = '  +3.37%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '3.119.85'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').Value = '''0.753'
$ws.Range('E11').Value = '  -2.27%  '
$ws.Range('D12').Value = '''0.203'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '''0.0000252'
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').Value = '''35.09'
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('D15').Value = '''5.59'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '91.578.21'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').Value = '3.715.65'
$ws.Range('E17').Value = '  +0.70%  '
$ws.Range('D18').Value = '3.146.20'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').Value = '''3.79'
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('D20').Value = '''14.87'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = '''5.82'
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('D22').Value = '''453.98'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E23').Value = '  -6.28%  '
$ws.Range('D24').Value = '''9.17'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '''5.88'
$ws.Range('E25').Value = '  -0.55%  '
$ws.Range('D26').Value = '''89.42'
$ws.Range('E26').Value = '  -4.06%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '''11.71'
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.291.32'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '''0.140'
$ws.Range('E30').Value = '  +16.42%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '''1.10'
$ws.Range('E31').Value = '  +9.60%  '
$ws.Range('E32').Value = '  -7.92%  '
$ws.Range('E33').Value = '  -8.61%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '''2.36'
$ws.Range('E34').Value = '  +22.09%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '''9.32'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.173'
$ws.Range('E36').Value = '  +3.39%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '''7.63'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '''26.32'
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('D39').Value = '''490.32'
$ws.Range('E39').Value = '  -1.55%  '
$ws.Range('B40').Value = 'MantraDAO'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D40').Value = '''3.90'
$ws.Range('E40').Value = '  -7.01%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '''1.31'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').Value = '''0.437'
$ws.Range('E42').Value = '  +3.13%  '
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '''157.69'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('D47').Value = '''1.92'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D48').Value = '''0.702'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('E49').Value = '  -1.09%  '
$ws.Range('D50').Value = '''4.45'
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('D51').Value = '''44.23'
$ws.Range('E51').Value = '  -1.33%  '
